# Weekly fruit/vegetable price update:
#  - Insert two new price records at the top of the data block (rows 570-571),
#    pushing all existing records for this sheet down by two rows.
#  - Append two new price records at the bottom of the data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows before row 570 (existing data shifts down) ---
$ws.Range("A570:A571").EntireRow.Insert()

# --- 2. Populate the two newly-inserted rows ---
$ws.Range("A570").Value = 3
$ws.Range("B570").Value = "Femacal de La Calera"
$ws.Range("C570").Value = "Coquimbo"
$ws.Range("D570").Value = 44753
$ws.Range("E570").Value = 5
$ws.Range("F570").Value = 100112045
$ws.Range("G570").Value = "Zapallo"
$ws.Range("H570").Value = "Camote"
$ws.Range("I570").Value = "1a (guarda)"
$ws.Range("J570").Value = 288
$ws.Range("K570").Value = 650
$ws.Range("L570").Value = 680
$ws.Range("M570").Value = 666
$ws.Range("N570").Value = "$/kilo (volumen en unidades)"
$ws.Range("O570").Value = "Provincia de Talca"
$ws.Range("P570").Value = 666
$ws.Range("Q570").Value = 1
$ws.Range("R570").Value = "Hortaliza"

$ws.Range("A571").Value = 3
$ws.Range("B571").Value = "Femacal de La Calera"
$ws.Range("C571").Value = "Coquimbo"
$ws.Range("D571").Value = 44753
$ws.Range("E571").Value = 5
$ws.Range("F571").Value = 100112045
$ws.Range("G571").Value = "Zapallo"
$ws.Range("H571").Value = "Paine"
$ws.Range("I571").Value = "1a (guarda)"
$ws.Range("J571").Value = 240
$ws.Range("K571").Value = 350
$ws.Range("L571").Value = 380
$ws.Range("M571").Value = 366
$ws.Range("N571").Value = "$/kilo (volumen en unidades)"
$ws.Range("O571").Value = "Provincia de Talca"
$ws.Range("P571").Value = 366
$ws.Range("Q571").Value = 1
$ws.Range("R571").Value = "Hortaliza"

# --- 3. Append two new rows (600 and 601) at the end of the data block ---
$ws.Range("A600").Value = 3
$ws.Range("B600").Value = "Femacal de La Calera"
$ws.Range("C600").Value = "Coquimbo"
$ws.Range("D600").Value = 44454
$ws.Range("E600").Value = 5
$ws.Range("F600").Value = 100112045
$ws.Range("G600").Value = "Zapallo"
$ws.Range("H600").Value = "Camote"
$ws.Range("I600").Value = "1a (guarda)"
$ws.Range("J600").Value = 120
$ws.Range("K600").Value = 800
$ws.Range("L600").Value = 800
$ws.Range("M600").Value = 800
$ws.Range("N600").Value = "$/kilo (volumen en unidades)"
$ws.Range("O600").Value = "Provincia de Talca"
$ws.Range("P600").Value = 800
$ws.Range("Q600").Value = 1
$ws.Range("R600").Value = "Hortaliza"

$ws.Range("A601").Value = 3
$ws.Range("B601").Value = "Femacal de La Calera"
$ws.Range("C601").Value = "Coquimbo"
$ws.Range("D601").Value = 44454
$ws.Range("E601").Value = 5
$ws.Range("F601").Value = 100112045
$ws.Range("G601").Value = "Zapallo"
$ws.Range("H601").Value = "Camote"
$ws.Range("I601").Value = "2a (guarda)"
$ws.Range("J601").Value = 100
$ws.Range("K601").Value = 600
$ws.Range("L601").Value = 600
$ws.Range("M601").Value = 600
$ws.Range("N601").Value = "$/kilo (volumen en unidades)"
$ws.Range("O601").Value = "Provincia de Talca"
$ws.Range("P601").Value = 600
$ws.Range("Q601").Value = 1
$ws.Range("R601").Value = "Hortaliza"
